$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 0.06722580394740553
$ws.Cells.Item(2, 8).Value = 3.821060765133717
$ws.Cells.Item(2, 9).Value = -70.19123444635315
$ws.Cells.Item(3, 7).Value = 0.07449243302246021
$ws.Cells.Item(3, 8).Value = 8.898991562842772
$ws.Cells.Item(4, 7).Value = -0.03892626210420437
$ws.Cells.Item(4, 8).Value = 13.29109288731735
$ws.Cells.Item(5, 7).Value = -0.04375896522075896
$ws.Cells.Item(5, 8).Value = -61.6401439174338
$ws.Cells.Item(6, 7).Value = -0.1059318550642394
$ws.Cells.Item(6, 8).Value = 0.1031139484681808
$ws.Cells.Item(7, 7).Value = -0.08580507744231732
$ws.Cells.Item(7, 8).Value = 6.091736149251741
$ws.Cells.Item(8, 7).Value = -0.3604675415178827
$ws.Cells.Item(8, 8).Value = 1.754404045009248
$ws.Cells.Item(9, 7).Value = -0.3819746142315359
$ws.Cells.Item(9, 8).Value = 2.0805699756689
$ws.Cells.Item(10, 7).Value = 0.02623607185697606
$ws.Cells.Item(10, 8).Value = 29.93945099451733
$ws.Cells.Item(11, 7).Value = 0.04085606379846735
$ws.Cells.Item(11, 8).Value = 80.02646863065915
$ws.Cells.Item(12, 7).Value = 0.218207108446705
$ws.Cells.Item(12, 8).Value = -1.596609521937906
$ws.Cells.Item(13, 7).Value = 0.2329176254027295
$ws.Cells.Item(13, 8).Value = 3.423885947408133
$ws.Cells.Item(14, 7).Value = -0.04046412012738942
$ws.Cells.Item(14, 8).Value = 3.897348487160356
$ws.Cells.Item(15, 7).Value = -0.04546777015092802
$ws.Cells.Item(15, 8).Value = 4.672325728382341
$ws.Cells.Item(16, 7).Value = 0.2219989499483603
$ws.Cells.Item(16, 8).Value = 4.434021956836862
$ws.Cells.Item(17, 7).Value = 0.2131858790883665
$ws.Cells.Item(17, 8).Value = -3.338088081131141
$ws.Cells.Item(18, 7).Value = 0.07713202805337856
$ws.Cells.Item(18, 8).Value = 5.631024259330906
$ws.Cells.Item(19, 7).Value = 0.07288161929369785
$ws.Cells.Item(19, 8).Value = -3.260496798340098
$ws.Cells.Item(20, 7).Value = -0.07321650954047218
$ws.Cells.Item(20, 8).Value = 2.354521260955603
$ws.Cells.Item(21, 7).Value = -0.08155183302143035
$ws.Cells.Item(21, 8).Value = 5.795386496671771
$ws.Cells.Item(22, 7).Value = 0.06720473908328339
$ws.Cells.Item(22, 8).Value = -8.567340466145298
$ws.Cells.Item(23, 7).Value = 0.07367378614657269
$ws.Cells.Item(23, 8).Value = 7.818479254017722
$ws.Cells.Item(24, 7).Value = 0.06398582350122863
$ws.Cells.Item(24, 8).Value = -3.94268677831709
$ws.Cells.Item(25, 7).Value = 0.06604577092699963
$ws.Cells.Item(25, 8).Value = 20.56536306401558
$ws.Cells.Item(26, 7).Value = 0.1210477642112683
$ws.Cells.Item(26, 8).Value = 1.423329481646041
$ws.Cells.Item(27, 7).Value = 0.128054879230127
$ws.Cells.Item(27, 8).Value = 12.46677769627858
$ws.Cells.Item(28, 7).Value = 0.1360149495921302
$ws.Cells.Item(28, 8).Value = 5.231252556964091
$ws.Cells.Item(29, 7).Value = 0.145122885947872
$ws.Cells.Item(29, 8).Value = -3.790447832968794
$ws.Cells.Item(30, 7).Value = 0.09081001026311288
$ws.Cells.Item(30, 8).Value = 7.713095788268538
$ws.Cells.Item(31, 7).Value = 0.09081001026311288
$ws.Cells.Item(31, 8).Value = 11.16683055072648
$ws.Cells.Item(32, 7).Value = 0.04965845441129226
$ws.Cells.Item(32, 8).Value = -6.936750040054932
$ws.Cells.Item(33, 7).Value = 0.05688078442867719
$ws.Cells.Item(33, 8).Value = 2.96485751582062
$ws.Cells.Item(34, 7).Value = 0.02245745151230089
$ws.Cells.Item(34, 8).Value = 29.38373951780778
$ws.Cells.Item(35, 7).Value = 0.02245745151230089
$ws.Cells.Item(35, 8).Value = 32.87891789432231
$ws.Cells.Item(36, 7).Value = -0.028389605615176
$ws.Cells.Item(36, 8).Value = 2.258836502252722
$ws.Cells.Item(37, 7).Value = -0.03407690386446319
$ws.Cells.Item(37, 8).Value = -2.444168337447866
$ws.Cells.Item(38, 7).Value = 0.07773205036937496
$ws.Cells.Item(38, 8).Value = -0.7034408144711588
$ws.Cells.Item(39, 7).Value = 0.07139458882708651
$ws.Cells.Item(39, 8).Value = -8.168411352366942
$ws.Cells.Item(40, 7).Value = 0.0671628500838474
$ws.Cells.Item(40, 8).Value = 1.440616835462994
$ws.Cells.Item(41, 7).Value = 0.07546912863120987
$ws.Cells.Item(41, 8).Value = 16.06556905782847
$ws.Cells.Item(42, 7).Value = 0.08158600782110707
$ws.Cells.Item(42, 8).Value = 4.878071234949394
$ws.Cells.Item(43, 7).Value = 0.08667070182049637
$ws.Cells.Item(43, 8).Value = 8.116074346633642
$ws.Cells.Item(44, 7).Value = 0.09059160601318694
$ws.Cells.Item(44, 8).Value = 2.656516249927409
$ws.Cells.Item(45, 7).Value = 0.08841057438685433
$ws.Cells.Item(45, 8).Value = -2.182844331727842
$ws.Cells.Item(46, 7).Value = 0.004174468531392874
$ws.Cells.Item(46, 8).Value = 252.5632568879498
$ws.Cells.Item(47, 7).Value = -0.00214766700859127
$ws.Cells.Item(47, 8).Value = -4395.529260849969
$ws.Cells.Item(48, 7).Value = -0.09354761524741659
$ws.Cells.Item(48, 8).Value = 2.661178059414909
$ws.Cells.Item(49, 7).Value = -0.0997141241738095
$ws.Cells.Item(49, 8).Value = 9.000201365396189
$ws.Cells.Item(50, 7).Value = 0.1640335734904538
$ws.Cells.Item(50, 8).Value = -3.792225421515707
$ws.Cells.Item(51, 7).Value = 0.1792571562026003
$ws.Cells.Item(51, 8).Value = 5.54791113014713
$ws.Cells.Item(52, 7).Value = 0.06177282845617584
$ws.Cells.Item(52, 8).Value = -12.94928423636981
$ws.Cells.Item(53, 7).Value = 0.06799973052480951
$ws.Cells.Item(53, 8).Value = 5.727371476420284
$ws.Cells.Item(54, 7).Value = -0.155716748650149
$ws.Cells.Item(54, 8).Value = -21.84136876878402
$ws.Cells.Item(55, 7).Value = -0.1274660012382811
$ws.Cells.Item(55, 8).Value = -9.439961732464022
$ws.Cells.Item(56, 7).Value = 0.1875307805177246
$ws.Cells.Item(56, 8).Value = -1.310411358668418
$ws.Cells.Item(57, 7).Value = 0.2028489696711937
$ws.Cells.Item(57, 8).Value = 1.981162082473394
